$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Sputum Test Ordered) updated values
$ws.Range("B3").Value = 56
$ws.Range("C3").Value = 250
$ws.Range("D3").Value = 0.224
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 179
$ws.Range("G3").Value = 0.15642458100558659
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 71
$ws.Range("J3").Value = 0.39436619718309862
$ws.Range("K3").Value = 3.7390578596597623
$ws.Range("L3").Value = 1.4685474569468147
$ws.Range("M3").Value = 3.3578726911263548
$ws.Range("N3").Value = 0.00078544779410171666
$ws.Range("O3").Value = 1.7315912590284546
$ws.Range("P3").Value = 8.0738185787144268

# Row 4 (Referral) updated values
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 250
$ws.Range("D4").Value = 0.056000000000000001
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 179
$ws.Range("G4").Value = 0.033519553072625698
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 71
$ws.Range("J4").Value = 0.1126760563380282
$ws.Range("K4").Value = 3.5945943348961187
$ws.Range("L4").Value = 2.0701756354350751
$ws.Range("M4").Value = 2.2215679969295099
$ws.Range("N4").Value = 0.026312514671228783
$ws.Range("O4").Value = 1.1626028168804754
$ws.Range("P4").Value = 11.1139490158277

# Row 6 (Antibiotics) updated values
$ws.Range("C6").Value = 250
$ws.Range("D6").Value = 0.043999999999999997
$ws.Range("F6").Value = 179
$ws.Range("G6").Value = 0.0167597765363128
$ws.Range("I6").Value = 71
$ws.Range("J6").Value = 0.1126760563380282
$ws.Range("K6").Value = 8.7196803624762946
$ws.Range("L6").Value = 6.3810709023567025
$ws.Range("M6").Value = 2.9592506020988121
$ws.Range("N6").Value = 0.003083882065132851
$ws.Range("O6").Value = 2.0777362438253051
$ws.Range("P6").Value = 36.594070036421584

$wb.Save()
